$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 6.6728885
$ws.Range("N2").Value = 13.345777
$ws.Range("O2").Value = 0.420360621011459
$ws.Range("P2").Value = 0.387781561837984
$ws.Range("Q2").Value = 0.1037478461018333
$ws.Range("R2").Value = 0.622487076611
$ws.Range("S2").Value = 0.420360621011459
$ws.Range("T2").Value = 0.387781561837984

# Row 3
$ws.Range("O3").Value = 0.07508775324202954
$ws.Range("P3").Value = 0.1039023808546958
$ws.Range("S3").Value = 0.07508775324202954
$ws.Range("T3").Value = 0.1039023808546958

# Row 4
$ws.Range("M4").Value = 0.7168596666666667
$ws.Range("N4").Value = 2.150579
$ws.Range("O4").Value = 0.04515879062838642
$ws.Range("P4").Value = 0.06248829749485322
$ws.Range("Q4").Value = 0.01114549514411111
$ws.Range("R4").Value = 0.100309456297
$ws.Range("S4").Value = 0.04515879062838642
$ws.Range("T4").Value = 0.06248829749485322

# Row 5
$ws.Range("M5").Value = 6.534003
$ws.Range("N5").Value = 13.068006
$ws.Range("O5").Value = 0.4116114871049826
$ws.Range("P5").Value = 0.3797105089338856
$ws.Range("Q5").Value = 0.101588500643
$ws.Range("R5").Value = 0.609531003858
$ws.Range("S5").Value = 0.4116114871049826
$ws.Range("T5").Value = 0.3797105089338856

# Row 6
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09023766666666666
$ws.Range("N6").Value = 0.270713
$ws.Range("O6").Value = 0.005684548992332936
$ws.Range("P6").Value = 0.00786597213109781
$ws.Range("Q6").Value = 0.001402985162111111
$ws.Range("R6").Value = 0.012626866459
$ws.Range("S6").Value = 0.005684548992332936
$ws.Range("T6").Value = 0.00786597213109781

# Row 7
$ws.Range("M7").Value = 0.668253
$ws.Range("N7").Value = 2.004759
$ws.Range("O7").Value = 0.04209679902080943
$ws.Range("P7").Value = 0.05825127874748356
$ws.Range("Q7").Value = 0.010389774893
$ws.Range("R7").Value = 0.09350797403699999
$ws.Range("S7").Value = 0.04209679902080943
$ws.Range("T7").Value = 0.05825127874748356
